# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets.
# Mapping of row -> new value for column F:
#   F2: 1371 -> 1373
#   F3: 2115 -> 2133
#   F4: 292  -> 298
#   F6: 6383 -> 6387
#   F7: 265  -> 269
#   F8: 117  -> 118

$wb = $excel.ActiveWorkbook

$updates = @{
    2 = 1373
    3 = 2133
    4 = 298
    6 = 6387
    7 = 269
    8 = 118
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
